$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.62046645117295
$ws.Range("C2").Value = 9.724763616072106
$ws.Range("E2").Value = 10.7878164373659
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 32.54262084097724
$ws.Range("H2").Value = 15.13093294753925
$ws.Range("I2").Value = 22.78502735584803
$ws.Range("L2").Value = 9.998072943991669
$ws.Range("B3").Value = 16.99069793147803
$ws.Range("C3").Value = 9.313849223107676
$ws.Range("E3").Value = 10.82561022885022
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 32.66164659694616
$ws.Range("H3").Value = 15.21542247173364
$ws.Range("I3").Value = 22.96726131102758
$ws.Range("L3").Value = 9.97096373439955
$ws.Range("B4").Value = 16.5940517925959
$ws.Range("C4").Value = 9.050361998562131
$ws.Range("E4").Value = 10.85038888027536
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 32.75352010666416
$ws.Range("H4").Value = 15.27163499439816
$ws.Range("I4").Value = 23.0865002819586
$ws.Range("L4").Value = 9.956190575717025
$ws.Range("B5").Value = 16.43015215960765
$ws.Range("C5").Value = 8.940267621025136
$ws.Range("E5").Value = 10.86088232704079
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 32.79563241182682
$ws.Range("H5").Value = 15.29562779856365
$ws.Range("I5").Value = 23.13693289393083
$ws.Range("L5").Value = 9.950645045410047
$ws.Range("B6").Value = 16.40280762275018
$ws.Range("C6").Value = 8.921825108746329
$ws.Range("E6").Value = 10.86264868112749
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 32.80290578096026
$ws.Range("H6").Value = 15.29967723150012
$ws.Range("I6").Value = 23.14541826968934
$ws.Range("L6").Value = 9.949752988064446
$ws.Range("B7").Value = 16.59185020991713
$ws.Range("C7").Value = 9.048888114964296
$ws.Range("E7").Value = 10.85052879489772
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 32.75406919831043
$ws.Range("H7").Value = 15.27195417955667
$ws.Range("I7").Value = 23.0871729838823
$ws.Range("L7").Value = 9.956113860058364
$ws.Range("B8").Value = 17.40552909537951
$ws.Range("C8").Value = 9.585452222204305
$ws.Range("E8").Value = 10.80052146196434
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 32.57973057961043
$ws.Range("H8").Value = 15.15916237604154
$ws.Range("I8").Value = 22.84633341500292
$ws.Range("L8").Value = 9.988339634568238
$ws.Range("B9").Value = 18.91258995164164
$ws.Range("C9").Value = 10.54555050614659
$ws.Range("E9").Value = 10.71492229519224
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 32.38909182041005
$ws.Range("H9").Value = 14.97258079905417
$ws.Range("I9").Value = 22.43259792142383
$ws.Range("L9").Value = 10.06619715723517
$ws.Range("B10").Value = 19.95472136127479
$ws.Range("C10").Value = 11.19107713815662
$ws.Range("E10").Value = 10.659608301623
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 32.34385013306392
$ws.Range("H10").Value = 14.85687219879648
$ws.Range("I10").Value = 22.16468729137451
$ws.Range("L10").Value = 10.13205649145619
$ws.Range("B11").Value = 20.41275152379852
$ws.Range("C11").Value = 11.47115639621009
$ws.Range("E11").Value = 10.63608475995331
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 32.34431946258776
$ws.Range("H11").Value = 14.80893486993093
$ws.Range("I11").Value = 22.0507231770008
$ws.Range("L11").Value = 10.1638290137541
$ws.Range("B12").Value = 20.58375428830663
$ws.Range("C12").Value = 11.57522433708122
$ws.Range("E12").Value = 10.62741238155755
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 32.34755420492404
$ws.Range("H12").Value = 14.79146309847263
$ws.Range("I12").Value = 22.00871312180251
$ws.Range("L12").Value = 10.17611468887578
$ws.Range("B13").Value = 20.54703644277129
$ws.Range("C13").Value = 11.55290060417304
$ws.Range("E13").Value = 10.62926966393389
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 32.34672116467065
$ws.Range("H13").Value = 14.79519558171655
$ws.Range("I13").Value = 22.01770964893696
$ws.Range("L13").Value = 10.17345755568009
$ws.Range("B14").Value = 20.42686971304789
$ws.Range("C14").Value = 11.47975825610323
$ws.Range("E14").Value = 10.63536656008985
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 32.34452416735288
$ws.Range("H14").Value = 14.80748377677509
$ws.Range("I14").Value = 22.047243974712
$ws.Range("L14").Value = 10.1648347097834
$ws.Range("B15").Value = 20.35294217017033
$ws.Range("C15").Value = 11.43469596746273
$ws.Range("E15").Value = 10.63913174464654
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 32.34357735358765
$ws.Range("H15").Value = 14.8150995032727
$ws.Range("I15").Value = 22.06548405493276
$ws.Range("L15").Value = 10.1595858593128
$ws.Range("B16").Value = 19.92445314671982
$ws.Range("C16").Value = 11.17249657131688
$ws.Range("E16").Value = 10.66117858363341
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 32.34424573682854
$ws.Range("H16").Value = 14.86010002861582
$ws.Range("I16").Value = 22.17229497768075
$ws.Range("L16").Value = 10.13001600738093
$ws.Range("B17").Value = 19.65737770896011
$ws.Range("C17").Value = 11.00813834369379
$ws.Range("E17").Value = 10.67512325065841
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 32.35007031328874
$ws.Range("H17").Value = 14.8889137251711
$ws.Range("I17").Value = 22.23985167187759
$ws.Range("L17").Value = 10.11233560426161
$ws.Range("B18").Value = 19.50226015470343
$ws.Range("C18").Value = 10.91232810195232
$ws.Range("E18").Value = 10.68329814399588
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 32.35540032737818
$ws.Range("H18").Value = 14.90592858124488
$ws.Range("I18").Value = 22.27945248524113
$ws.Range("L18").Value = 10.10233734858269
$ws.Range("B19").Value = 19.44948635954511
$ws.Range("C19").Value = 10.8796706273279
$ws.Range("E19").Value = 10.68609252868992
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 32.35754401689267
$ws.Range("H19").Value = 14.91176526773633
$ws.Range("I19").Value = 22.29298817495079
$ws.Range("L19").Value = 10.09898168204275
$ws.Range("B20").Value = 19.68596492798788
$ws.Range("C20").Value = 11.02576690354459
$ws.Range("E20").Value = 10.6736228503014
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 32.34924513301841
$ws.Range("H20").Value = 14.885800676389
$ws.Range("I20").Value = 22.23258308783923
$ws.Range("L20").Value = 10.11420005730109
$ws.Range("B21").Value = 20.46223290857537
$ws.Range("C21").Value = 11.50129628434713
$ws.Range("E21").Value = 10.63356936483988
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 32.34508630036451
$ws.Range("H21").Value = 14.80385590973141
$ws.Range("I21").Value = 22.03853786620575
$ws.Range("L21").Value = 10.16736060628126
$ws.Range("B22").Value = 20.95527592652889
$ws.Range("C22").Value = 11.80045740521205
$ws.Range("E22").Value = 10.60876453833135
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 32.36019684228493
$ws.Range("H22").Value = 14.75427247203437
$ws.Range("I22").Value = 21.91840009405692
$ws.Range("L22").Value = 10.20358195195574
$ws.Range("B23").Value = 20.69347813397865
$ws.Range("C23").Value = 11.64186466619555
$ws.Range("E23").Value = 10.62187783249281
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 32.35049227339657
$ws.Range("H23").Value = 14.78037085905793
$ws.Range("I23").Value = 21.98190565779669
$ws.Range("L23").Value = 10.18411697595592
$ws.Range("B24").Value = 19.67304553315482
$ws.Range("C24").Value = 11.01780113544689
$ws.Range("E24").Value = 10.67430068939624
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 32.34961202758838
$ws.Range("H24").Value = 14.88720668604995
$ws.Range("I24").Value = 22.23586684306413
$ws.Range("L24").Value = 10.1133566191463
$ws.Range("B25").Value = 18.51560330636324
$ws.Range("C25").Value = 10.29610104025474
$ws.Range("E25").Value = 10.73674711813148
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 32.42418513209
$ws.Range("H25").Value = 15.01932304845897
$ws.Range("I25").Value = 22.53821940065141
$ws.Range("L25").Value = 10.04359225348672
